$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema3c"
$ws.Cells.Item(2,3).Value = "Nrp2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.2364093333333333
$ws.Cells.Item(2,8).Value = 0.709228
$ws.Cells.Item(2,9).Value = 0.005805733041453686
$ws.Cells.Item(2,10).Value = 0.005805733041453687
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 34.52052333333333
$ws.Cells.Item(2,14).Value = 103.56157
$ws.Cells.Item(2,15).Value = 0.7684334662422598
$ws.Cells.Item(2,16).Value = 0.7684334662422598
$ws.Cells.Item(2,17).Value = 8.16097390755111
$ws.Cells.Item(2,18).Value = 73.44876516795999
$ws.Cells.Item(2,19).Value = 0.004461319565121474
$ws.Cells.Item(2,20).Value = 0.004461319565121475

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema3c"
$ws.Cells.Item(3,3).Value = "Nrp2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.2364093333333333
$ws.Cells.Item(3,8).Value = 0.709228
$ws.Cells.Item(3,9).Value = 0.005805733041453686
$ws.Cells.Item(3,10).Value = 0.005805733041453687
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 5.347618333333333
$ws.Cells.Item(3,14).Value = 16.042855
$ws.Cells.Item(3,15).Value = 0.1190390091234806
$ws.Cells.Item(3,16).Value = 0.1190390091234805
$ws.Cells.Item(3,17).Value = 1.264226885104444
$ws.Cells.Item(3,18).Value = 11.37804196594
$ws.Cells.Item(3,19).Value = 0.0006911087084900979
$ws.Cells.Item(3,20).Value = 0.0006911087084900979

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema3c"
$ws.Cells.Item(4,3).Value = "Nrp2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.2364093333333333
$ws.Cells.Item(4,8).Value = 0.709228
$ws.Cells.Item(4,9).Value = 0.005805733041453686
$ws.Cells.Item(4,10).Value = 0.005805733041453687
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.055101333333334
$ws.Cells.Item(4,14).Value = 15.165304
$ws.Cells.Item(4,15).Value = 0.1125275246342597
$ws.Cells.Item(4,16).Value = 0.1125275246342597
$ws.Cells.Item(4,17).Value = 1.195073136145778
$ws.Cells.Item(4,18).Value = 10.755658225312
$ws.Cells.Item(4,19).Value = 0.0006533047678421153
$ws.Cells.Item(4,20).Value = 0.0006533047678421153

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Sema3c"
$ws.Cells.Item(5,3).Value = "Nrp2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 39.09670133333334
$ws.Cells.Item(5,8).Value = 117.290104
$ws.Cells.Item(5,9).Value = 0.9601355730855794
$ws.Cells.Item(5,10).Value = 0.9601355730855795
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 34.52052333333333
$ws.Cells.Item(5,14).Value = 103.56157
$ws.Cells.Item(5,15).Value = 0.7684334662422598
$ws.Cells.Item(5,16).Value = 0.7684334662422598
$ws.Cells.Item(5,17).Value = 1349.638590633698
$ws.Cells.Item(5,18).Value = 12146.74731570328
$ws.Cells.Item(5,19).Value = 0.7378003064886504
$ws.Cells.Item(5,20).Value = 0.7378003064886505

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sema3c"
$ws.Cells.Item(6,3).Value = "Nrp2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 39.09670133333334
$ws.Cells.Item(6,8).Value = 117.290104
$ws.Cells.Item(6,9).Value = 0.9601355730855794
$ws.Cells.Item(6,10).Value = 0.9601355730855795
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 5.347618333333333
$ws.Cells.Item(6,14).Value = 16.042855
$ws.Cells.Item(6,15).Value = 0.1190390091234806
$ws.Cells.Item(6,16).Value = 0.1190390091234805
$ws.Cells.Item(6,17).Value = 209.0742368229911
$ws.Cells.Item(6,18).Value = 1881.66813140692
$ws.Cells.Item(6,19).Value = 0.1142935872443125
$ws.Cells.Item(6,20).Value = 0.1142935872443125

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema3c"
$ws.Cells.Item(7,3).Value = "Nrp2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 39.09670133333334
$ws.Cells.Item(7,8).Value = 117.290104
$ws.Cells.Item(7,9).Value = 0.9601355730855794
$ws.Cells.Item(7,10).Value = 0.9601355730855795
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.055101333333334
$ws.Cells.Item(7,14).Value = 15.165304
$ws.Cells.Item(7,15).Value = 0.1125275246342597
$ws.Cells.Item(7,16).Value = 0.1125275246342597
$ws.Cells.Item(7,17).Value = 197.6377870390685
$ws.Cells.Item(7,18).Value = 1778.740083351616
$ws.Cells.Item(7,19).Value = 0.1080416793526166
$ws.Cells.Item(7,20).Value = 0.1080416793526166

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Sema3c"
$ws.Cells.Item(8,3).Value = "Nrp2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.386869333333333
$ws.Cells.Item(8,8).Value = 4.160608
$ws.Cells.Item(8,9).Value = 0.03405869387296686
$ws.Cells.Item(8,10).Value = 0.03405869387296687
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 34.52052333333333
$ws.Cells.Item(8,14).Value = 103.56157
$ws.Cells.Item(8,15).Value = 0.7684334662422598
$ws.Cells.Item(8,16).Value = 0.7684334662422598
$ws.Cells.Item(8,17).Value = 47.87545518161777
$ws.Cells.Item(8,18).Value = 430.8790966345599
$ws.Cells.Item(8,19).Value = 0.02617184018848794
$ws.Cells.Item(8,20).Value = 0.02617184018848795

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Sema3c"
$ws.Cells.Item(9,3).Value = "Nrp2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.386869333333333
$ws.Cells.Item(9,8).Value = 4.160608
$ws.Cells.Item(9,9).Value = 0.03405869387296686
$ws.Cells.Item(9,10).Value = 0.03405869387296687
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.347618333333333
$ws.Cells.Item(9,14).Value = 16.042855
$ws.Cells.Item(9,15).Value = 0.1190390091234806
$ws.Cells.Item(9,16).Value = 0.1190390091234805
$ws.Cells.Item(9,17).Value = 7.41644787287111
$ws.Cells.Item(9,18).Value = 66.74803085584
$ws.Cells.Item(9,19).Value = 0.004054313170677933
$ws.Cells.Item(9,20).Value = 0.004054313170677934

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Sema3c"
$ws.Cells.Item(10,3).Value = "Nrp2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.386869333333333
$ws.Cells.Item(10,8).Value = 4.160608
$ws.Cells.Item(10,9).Value = 0.03405869387296686
$ws.Cells.Item(10,10).Value = 0.03405869387296687
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.055101333333334
$ws.Cells.Item(10,14).Value = 15.165304
$ws.Cells.Item(10,15).Value = 0.1125275246342597
$ws.Cells.Item(10,16).Value = 0.1125275246342597
$ws.Cells.Item(10,17).Value = 7.010765016092445
$ws.Cells.Item(10,18).Value = 63.096885144832
$ws.Cells.Item(10,19).Value = 0.003832540513800988
$ws.Cells.Item(10,20).Value = 0.003832540513800988
